$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.306.59"
$ws.Range("E2").Value = "  +1.32%  "
$ws.Range("D3").Value = "1.612.06"
$ws.Range("E3").Value = "  +0.83%  "
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "'213.24"
$ws.Range("E5").Value = "  +0.51%  "
$ws.Range("E6").Value = "  -0.24%  "
$ws.Range("E7").Value = "  +0.48%  "
$ws.Range("E8").Value = "  +1.18%  "
$ws.Range("E9").Value = "  +0.35%  "
$ws.Range("D10").Value = "'18.44"
$ws.Range("E10").Value = "  +2.36%  "
$ws.Range("E11").Value = "  +0.03%  "
$ws.Range("D12").Value = "1.835.88"
$ws.Range("E12").Value = "  +0.78%  "
$ws.Range("D13").Value = "1.631.63"
$ws.Range("E13").Value = "  +1.96%  "
$ws.Range("E14").Value = "  +0.32%  "
$ws.Range("E15").Value = "  +1.33%  "
$ws.Range("D16").Value = "26.290.33"
$ws.Range("E16").Value = "  +1.18%  "
$ws.Range("D17").Value = "'62.24"
$ws.Range("E17").Value = "  +3.39%  "
$ws.Range("D18").Value = "0.0₃0727"
$ws.Range("E18").Value = "  +1.08%  "
$ws.Range("E19").Value = "  -0.17%  "
$ws.Range("D20").Value = "'201.65"
$ws.Range("E20").Value = "  +0.77%  "
$ws.Range("E21").Value = "  +1.32%  "
$ws.Range("D22").Value = "'9.35"
$ws.Range("E22").Value = "  +1.01%  "
$ws.Range("E23").Value = "  +0.82%  "
$ws.Range("E24").Value = "  +1.93%  "
$ws.Range("D25").Value = "'143.52"
$ws.Range("E25").Value = "  +1.66%  "
$ws.Range("E26").Value = "  -0.19%  "
$ws.Range("E27").Value = "  -0.28%  "
$ws.Range("D28").Value = "'15.24"
$ws.Range("E28").Value = "  +0.69%  "
$ws.Range("D29").Value = "'6.56"
$ws.Range("E29").Value = "  +2.50%  "
$ws.Range("D30").Value = "'0.0499"
$ws.Range("E30").Value = "  +5.70%  "
$ws.Range("E31").Value = "  +0.38%  "
$ws.Range("E32").Value = "  +2.78%  "
$ws.Range("D33").Value = "'2.95"
$ws.Range("E33").Value = "  +0.07%  "
$ws.Range("E34").Value = "  +1.23%  "
$ws.Range("D35").Value = "'2.37"
$ws.Range("E35").Value = "  +0.79%  "
$ws.Range("D36").Value = "1.162.79"
$ws.Range("E36").Value = "  +3.46%  "
$ws.Range("E37").Value = "  +0.64%  "
$ws.Range("E38").Value = "  -0.19%  "
$ws.Range("E39").Value = "  +1.28%  "
$ws.Range("E40").Value = "  +0.29%  "
$ws.Range("E41").Value = "  +1.51%  "
$ws.Range("E42").Value = "  +4.25%  "
$ws.Range("E43").Value = "  +0.39%  "
$ws.Range("D44").Value = "1.746.53"
$ws.Range("E44").Value = "  +0.64%  "
$ws.Range("D45").Value = "'92.75"
$ws.Range("E45").Value = "  -0.32%  "
$ws.Range("E46").Value = "  +11.87%  "
$ws.Range("E47").Value = "  +1.84%  "
$ws.Range("D48").Value = "'53.87"
$ws.Range("E48").Value = "  +1.29%  "
$ws.Range("E49").Value = "  +0.74%  "
$ws.Range("E50").Value = "  -0.14%  "
$ws.Range("E51").Value = "  -0.29%  "
